$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking text (prices). Force it to remain
# text so values like "1.00", "0.999" or "0.0000225" keep their exact
# displayed form instead of Excel auto-converting them to numbers.
$ws.Range("D2:D51").NumberFormat = "@"

# Apply the updated coin price / volume(1h) data

# Row 2
$ws.Range("D2").Value = '69.505.65'
$ws.Range("E2").Value = '  +3.87%  '

# Row 3
$ws.Range("D3").Value = '3.640.55'
$ws.Range("E3").Value = '  +3.20%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.18%  '

# Row 5
$ws.Range("D5").Value = '627.71'
$ws.Range("E5").Value = '  +3.62%  '

# Row 6
$ws.Range("D6").Value = '160.22'
$ws.Range("E6").Value = '  +4.60%  '

# Row 7
$ws.Range("D7").Value = '3.639.84'
$ws.Range("E7").Value = '  +3.24%  '

# Row 8
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.19%  '

# Row 9
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  +2.53%  '

# Row 10
$ws.Range("D10").Value = '0.145'
$ws.Range("E10").Value = '  +3.99%  '

# Row 11
$ws.Range("D11").Value = '7.25'
$ws.Range("E11").Value = '  +6.03%  '

# Row 12
$ws.Range("D12").Value = '0.440'
$ws.Range("E12").Value = '  +3.25%  '

# Row 13
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").Value = '  +1.43%  '

# Row 14
$ws.Range("D14").Value = '33.46'
$ws.Range("E14").Value = '  +5.70%  '

# Row 15
$ws.Range("D15").Value = '4.252.66'
$ws.Range("E15").Value = '  +3.09%  '

# Row 16
$ws.Range("D16").Value = '3.638.60'
$ws.Range("E16").Value = '  +3.29%  '

# Row 17
$ws.Range("D17").Value = '69.298.47'
$ws.Range("E17").Value = '  +3.49%  '

# Row 18
$ws.Range("E18").Value = '  -0.16%  '

# Row 19
$ws.Range("D19").Value = '6.61'
$ws.Range("E19").Value = '  +5.47%  '

# Row 20
$ws.Range("D20").Value = '15.92'
$ws.Range("E20").Value = '  +4.05%  '

# Row 21
$ws.Range("D21").Value = '10.28'
$ws.Range("E21").Value = '  +11.10%  '

# Row 22
$ws.Range("D22").Value = '462.92'
$ws.Range("E22").Value = '  +4.17%  '

# Row 23
$ws.Range("D23").Value = '0.644'
$ws.Range("E23").Value = '  +2.85%  '

# Row 24
$ws.Range("B24").Value = 'PEPE'
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D24").Value = '0.0000138'
$ws.Range("E24").Value = '  +12.63%  '

# Row 25
$ws.Range("B25").Value = 'Litecoin'
$ws.Range("C25").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D25").Value = '78.67'
$ws.Range("E25").Value = '  +0.84%  '

# Row 26
$ws.Range("D26").Value = '3.781.24'
$ws.Range("E26").Value = '  +2.99%  '

# Row 27
$ws.Range("D27").Value = '10.57'
$ws.Range("E27").Value = '  +4.02%  '

# Row 28
$ws.Range("E28").Value = '  +0.19%  '

# Row 29
$ws.Range("D29").Value = '9.25'
$ws.Range("E29").Value = '  +13.60%  '

# Row 30
$ws.Range("B30").Value = 'Fetch.AI'
$ws.Range("C30").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D30").Value = '1.74'
$ws.Range("E30").Value = '  +5.76%  '

# Row 31
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").Value = '2.62'
$ws.Range("E31").Value = '  +3.80%  '

# Row 32
$ws.Range("D32").Value = '0.175'
$ws.Range("E32").Value = '  +10.88%  '

# Row 33
$ws.Range("D33").Value = '6.65'
$ws.Range("E33").Value = '  +8.57%  '

# Row 34
$ws.Range("D34").Value = '1.00'
$ws.Range("E34").Value = '  -0.14%  '

# Row 35
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Value = '1.97'
$ws.Range("E35").Value = '  +5.27%  '

# Row 36
$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '26.53'
$ws.Range("E36").Value = '  +3.61%  '

# Row 37
$ws.Range("D37").Value = '3.626.77'
$ws.Range("E37").Value = '  +2.95%  '

# Row 38
$ws.Range("D38").Value = '8.40'
$ws.Range("E38").Value = '  +5.38%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").Value = '2.38'
$ws.Range("E39").Value = '  +11.09%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D40").Value = '1.00'
$ws.Range("E40").Value = '  +0.02%  '

# Row 41
$ws.Range("D41").Value = '0.0932'
$ws.Range("E41").Value = '  +8.43%  '

# Row 42
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.27%  '

# Row 43
$ws.Range("D43").Value = '175.07'
$ws.Range("E43").Value = '  +0.02%  '

# Row 44
$ws.Range("D44").Value = '5.68'
$ws.Range("E44").Value = '  +2.83%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '31.33'
$ws.Range("E45").Value = '  +14.04%  '

# Row 46
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").Value = '0.917'
$ws.Range("E46").Value = '  +3.14%  '

# Row 47
$ws.Range("D47").Value = '1.38'
$ws.Range("E47").Value = '  +12.91%  '

# Row 48
$ws.Range("B48").Value = 'dogwifhat'
$ws.Range("C48").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D48").Value = '2.84'
$ws.Range("E48").Value = '  +9.32%  '

# Row 49
$ws.Range("B49").Value = 'OKB'
$ws.Range("C49").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D49").Value = '46.41'
$ws.Range("E49").Value = '  +1.63%  '

# Row 50
$ws.Range("B50").Value = 'TheGraph'
$ws.Range("C50").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D50").Value = '0.271'
$ws.Range("E50").Value = '  +9.56%  '

# Row 51
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D51").Value = '7.81'
$ws.Range("E51").Value = '  +3.55%  '
